$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RVC")

# ---------------------------------------------------------------------------
# The H:K "DY Xfmr zero-sequence" block is being shifted down two rows (to
# make room for a new "Grounding Calculations" heading at H2) and three of
# its formulas are re-pointed from a single shared $D4-anchored pattern to
# per-row D4/D5/D6 references. Clear the whole existing block first so the
# shared-formula remnants don't linger, then rewrite every cell at its new
# location.
# ---------------------------------------------------------------------------
$ws.Range("H1:K9").Clear()

# New section heading
$ws.Range("H2").Value = "Grounding Calculations"
$ws.Range("H2").Font.Bold = $true
$ws.Range("H2").Font.Italic = $true

# Row 3 (was row 1) - column headers
$ws.Range("J3").Value = "R0 [W]"
$ws.Range("J3").Font.Bold = $true
$ws.Range("J3").HorizontalAlignment = -4152
$ws.Range("K3").Value = "X0 [W]"
$ws.Range("K3").Font.Bold = $true
$ws.Range("K3").HorizontalAlignment = -4152

# Row 4 (was row 2)
$ws.Range("I4").Value = "DY Xfmr"
$ws.Range("I4").Font.Bold = $true
$ws.Range("I4").HorizontalAlignment = -4152
$ws.Range("J4").Formula = "=E2"
$ws.Range("J4").NumberFormat = "0.0000"
$ws.Range("K4").Formula = "=F2"
$ws.Range("K4").NumberFormat = "0.0000"

# Row 5 (was row 3)
$ws.Range("H5").Value = "R0/mi"
$ws.Range("H5").Font.Bold = $true
$ws.Range("I5").Value = "X0/mi"
$ws.Range("I5").Font.Bold = $true

# Row 6 (was row 4)
$ws.Range("H6").Value = 0.632
$ws.Range("H6").NumberFormat = "0.0000"
$ws.Range("I6").Value = 1.744
$ws.Range("I6").NumberFormat = "0.0000"
$ws.Range("J6").Formula = "=H6*`$D4/5280"
$ws.Range("J6").NumberFormat = "0.0000"
$ws.Range("K6").Formula = "=I6*`$D4/5280"
$ws.Range("K6").NumberFormat = "0.0000"

# Row 7 (was row 5)
$ws.Range("H7").Value = 2.392
$ws.Range("H7").NumberFormat = "0.0000"
$ws.Range("I7").Value = 2.568
$ws.Range("I7").NumberFormat = "0.0000"
$ws.Range("J7").Formula = "=H7*`$D5/5280"
$ws.Range("J7").NumberFormat = "0.0000"
$ws.Range("K7").Formula = "=I7*`$D5/5280"
$ws.Range("K7").NumberFormat = "0.0000"

# Row 8 (was row 6)
$ws.Range("H8").Value = 5.9947
$ws.Range("H8").NumberFormat = "0.0000"
$ws.Range("I8").Value = 3.0253
$ws.Range("I8").NumberFormat = "0.0000"
$ws.Range("J8").Formula = "=H8*`$D6/5280"
$ws.Range("J8").NumberFormat = "0.0000"
$ws.Range("K8").Formula = "=I8*`$D6/5280"
$ws.Range("K8").NumberFormat = "0.0000"

# Row 9 (was row 7)
$ws.Range("I9").Value = "Z0 pcc [W]"
$ws.Range("I9").Font.Bold = $true
$ws.Range("J9").Formula = "=SUM(J4:J8)"
$ws.Range("J9").NumberFormat = "0.0000"
$ws.Range("K9").Formula = "=SUM(K4:K8)"
$ws.Range("K9").NumberFormat = "0.0000"

# Row 10 (was row 8)
$ws.Range("H10").Value = "Check X0/X1 < 3?"
$ws.Range("J10").Formula = "=K9/F7"
$ws.Range("J10").NumberFormat = "0.00"
$ws.Range("K10").Value = "(pass)"

# Row 11 (was row 9)
$ws.Range("H11").Value = "Check R0/X1 < 1?"
$ws.Range("J11").Formula = "=J9/F7"
$ws.Range("J11").NumberFormat = "0.00"
$ws.Range("K11").Value = "(fail)"

# ---------------------------------------------------------------------------
# Column widths for the newly visible/used columns around the block
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 1.6
$ws.Columns.Item(9).ColumnWidth = 8.6
$ws.Columns.Item(10).ColumnWidth = 6.1
$ws.Columns.Item(11).ColumnWidth = 6.25

# Selection, matching where the author's screen shot selection ended up
$ws.Range("N28").Select()
